$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the cookie products and their picture filenames so the
# referenced images actually exist ("pictures can show now").
$ws.Range("B2").Value = "經典可可"
$ws.Range("C2").Value = "classic_cocoa.png"
$ws.Range("B3").Value = "經典抹茶"
$ws.Range("C3").Value = "classic_matcha.png"
$ws.Range("B4").Value = "愛戀玫瑰"
$ws.Range("C4").Value = "rose_love.png"

# Row 1 no longer has an explicit (taller) row height - let it size
# back to the sheet's default.
$ws.Rows("1").AutoFit()

# Move the active selection to C4.
$ws.Range("C4").Select()
